# Sprint_Burndown_4.xlsx edit script
# Updates the "Actual User Story Points Burned" log with additional entries
# (code review entries), extends the burndown formulas to cover the new
# row range, and updates the chart/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the cascading "Actual" burndown formulas (row 15) ---
# Previously summed K25:K30 / L25:L30 / ... ; now the log runs to row 37.
$ws.Range("D15").Formula = "=C15-SUM(K25:K37)"
$ws.Range("E15").Formula = "=D15-SUM(L25:L37)"
$ws.Range("F15").Formula = "=E15-SUM(M25:M37)"
$ws.Range("G15").Formula = "=F15-SUM(N25:N37)"
$ws.Range("H15").Formula = "=G15-SUM(O25:O37)"

# --- Clear the old "Actual User Story Points Burned" log rows (25-31) ---
# so the previous row 27 ("22b") / row 31 ("Total:") entries don't linger
# once the table is rebuilt with more rows below.
$ws.Range("J25:R31").ClearContents()

# --- Rebuild the log with the new/expanded set of entries (rows 25-37) ---
$ws.Range("J25").Value = 10
$ws.Range("O25").Value = 3

$ws.Range("J26").Value = 20
$ws.Range("M26").Value = 3

$ws.Range("J27").Value = 21
$ws.Range("M27").Value = 1

$ws.Range("J28").Value = "22b"
$ws.Range("M28").Value = 3

$ws.Range("J29").Value = "22c"
$ws.Range("O29").Value = 5

$ws.Range("J30").Value = 25
$ws.Range("K30").Value = 1

$ws.Range("J31").Value = 26
$ws.Range("K31").Value = 1

$ws.Range("J32").Value = 27
$ws.Range("M32").Value = 1

$ws.Range("J33").Value = 28
$ws.Range("M33").Value = 1

$ws.Range("J34").Value = 29
$ws.Range("M34").Value = 1

$ws.Range("J35").Value = 30
$ws.Range("M35").Value = 1

$ws.Range("J36").Value = 31
$ws.Range("M36").Value = 1

$ws.Range("J37").Value = 32
$ws.Range("M37").Value = 1

$ws.Range("J38").Value = "Total:"
$ws.Range("K38").Formula = "=SUM(K25:O37)"

# --- Update the selection shown when the workbook was last saved ---
$ws.Range("B3:G3").Select()
